$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# The legacy GSC export had a duplicated/blank leading day (2025-09-27) before
# the real data started on 2025-09-28. Remove that stray row so the whole
# table shifts up by one day; Excel re-keys the remaining rows (and the
# shared-string table) automatically.
$ws.Rows("2:2").Delete()

# Row 90 (now the last data row, 2025-12-25) inherited a legacy bad cell in
# the old export (an empty-string "Impressions" value instead of a numeric
# 0). Normalize it back to the numeric 0 used by every other row.
$ws.Cells.Item(90, 4).Value = 0
